$p = $ppt.ActivePresentation

# Slide 5 notes ("Parametre?") - clear the notes text body
$s5 = $p.Slides.Item(5)
$np5 = $s5.NotesPage
$np5.Shapes.Item(2).TextFrame.TextRange.Text = ""

# Slide 6 notes ("Mistral") - clear the notes text body
$s6 = $p.Slides.Item(6)
$np6 = $s6.NotesPage
$np6.Shapes.Item(2).TextFrame.TextRange.Text = ""
